# Refresh the cryptos "Price" (D) and "Volume(1h)" (E) columns with the
# latest scraped figures (scheduled GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.865.75"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "1.904.06"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").Value = "'313.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "'0.5033"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.34%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.07302"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'0.9099"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.37%  "

$ws.Range("D11").Value = "'20.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "'0.07660"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").Value = "1.884.52"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").Value = "'5.485"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "'6.612"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").Value = "'91.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'0.000008704"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").Value = "27.888.04"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").Value = "'14.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "

$ws.Range("D22").Value = "'5.161"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'10.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").Value = "'154.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("D25").Value = "'1.861"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.22%  "

$ws.Range("D26").Value = "'2.226"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.12%  "

$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("D28").Value = "'115.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("D29").Value = "'4.922"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "'0.08977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("E31").Value = "  -4.06%  "

$ws.Range("D32").Value = "'1.242"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "'0.7711"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("D34").Value = "'4.641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  +0.69%  "

$ws.Range("D36").Value = "'2.561"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").Value = "'1.100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").Value = "'0.5546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("D39").Value = "'3.017"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("D40").Value = "'0.05262"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").Value = "'6.968"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "

$ws.Range("D42").Value = "'8.535"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("D43").Value = "'0.1523"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").Value = "'111.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("D45").Value = "'10.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").Value = "'0.4799"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").Value = "'1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "'1.639"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("D49").Value = "'67.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("D50").Value = "'0.06086"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").Value = "'0.9011"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
